# Generate Report for Handback
# Appends a new row (for file 4a9157be-11bc-49be-a3e8-e9f6f83e2d16.md) to the
# "Overview", "zh-cn" and "de-de" tables/sheets of the handback-status report.

$wb = $excel.ActiveWorkbook

$fileGuid   = "4a9157be-11bc-49be-a3e8-e9f6f83e2d16"
$mdName     = "$fileGuid.md"
$mdPath     = "e2e\$fileGuid.md"
$zhXlfHash  = "fea957e22beea9dbd92753946ecd094df297fc44"
$zhXlfName  = "$fileGuid.$zhXlfHash.zh-cn.xlf"
$deXlfName  = "$fileGuid.$zhXlfHash.de-de.xlf"
$status     = "Handed back: in sync with en-US"

$dateHandoff   = "2016-08-12 14:53:26"
$dateZhHO      = "2016-08-12 14:53:19"
$dateZhHB      = "2016-08-12 14:53:47"
$dateDeHB      = "2016-08-12 14:53:55"

$oltestHash = "7bff609367b198cc453927766c9f5738d6da59a"
$zhcnHash   = "9b7413b55a8a1669d14f0229ab719caed3d88eb"
$dedeHash   = "a4a0ba8ff5abdc7063a615417d07d3de3aef9d4"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(4,1).Value = $mdName
$wsOverview.Cells.Item(4,2).Value = $mdPath
$wsOverview.Cells.Item(4,3).Value = ".md"
$wsOverview.Cells.Item(4,5).Value = $status
$wsOverview.Cells.Item(4,6).Value = $status
$wsOverview.Cells.Item(4,7).Value = $dateHandoff
$wsOverview.Cells.Item(4,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4,2), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$oltestHash/e2e/$mdName", "", "", $mdPath) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Cells.Item(4,1).Value  = $mdName
$wsZh.Cells.Item(4,2).Value  = ".md"
$wsZh.Cells.Item(4,3).Value  = $status
$wsZh.Cells.Item(4,4).Value  = "e2e"
$wsZh.Cells.Item(4,5).Value  = "ht"
$wsZh.Cells.Item(4,6).Value  = "True"
$wsZh.Cells.Item(4,7).Value  = $zhXlfName
$wsZh.Cells.Item(4,8).Value  = $dateZhHO
$wsZh.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,9).Value  = $mdName
$wsZh.Cells.Item(4,10).Value = $zhXlfName
$wsZh.Cells.Item(4,11).Value = $dateZhHB
$wsZh.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,12).Value = "'"
$wsZh.Cells.Item(4,13).Value = "True"
$wsZh.Cells.Item(4,14).Value = "'"
$wsZh.Cells.Item(4,15).Value = "False"
$wsZh.Cells.Item(4,16).Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$oltestHash/e2e/$mdName", "", "", $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4,9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$zhcnHash/e2e/$mdName", "", "", $mdName) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Cells.Item(4,1).Value  = $mdName
$wsDe.Cells.Item(4,2).Value  = ".md"
$wsDe.Cells.Item(4,3).Value  = $status
$wsDe.Cells.Item(4,4).Value  = "e2e"
$wsDe.Cells.Item(4,5).Value  = "ht"
$wsDe.Cells.Item(4,6).Value  = "True"
$wsDe.Cells.Item(4,7).Value  = $deXlfName
$wsDe.Cells.Item(4,8).Value  = $dateHandoff
$wsDe.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,9).Value  = $mdName
$wsDe.Cells.Item(4,10).Value = $deXlfName
$wsDe.Cells.Item(4,11).Value = $dateDeHB
$wsDe.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,12).Value = "'"
$wsDe.Cells.Item(4,13).Value = "True"
$wsDe.Cells.Item(4,14).Value = "'"
$wsDe.Cells.Item(4,15).Value = "False"
$wsDe.Cells.Item(4,16).Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4,1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$oltestHash/e2e/$mdName", "", "", $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4,9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$dedeHash/e2e/$mdName", "", "", $mdName) | Out-Null

Write-Output "Handback report row appended for $mdName"
